$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the last-updated date (C1) ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45387

# --- "BAU Emissions" sheet: refresh from the 4.0 model re-run ---
$bau = $wb.Worksheets.Item("BAU Emissions")

# Strip the old "NoSettings" suffix from every row label -> "test"
# (updates the shared-string labels in column A, rows 4:280)
[void]$bau.Cells.Replace(" : NoSettings", " : test")

# Updated values for row 94 (Industrial Sector Energy Related Emissions
# before CCS[natural gas if,iron and steel 241,CO2]) from the re-run
$cols = @("M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE")
$newVals = @(1001080, 2002150, 3003230, 4004300, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380, 5005380)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $bau.Range($cols[$i] + "94").Value = $newVals[$i]
}

# Scroll/select down near the bottom of the refreshed data
$bau.Activate()
[void]$bau.Range("A30:AE280").Select()

# --- "Current and Planned Capacity" loses the active tab, "About" gains it ---
$about.Activate()
